# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table (rows 2-50) with the latest coinranking.com snapshot, as produced
# by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.995.42'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.827.87'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("D4").Value = '''0.9959'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = '''243.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '''0.6320'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '''0.9986'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '''0.07501'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = '''0.2942'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = '''23.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").Value = '''0.07702'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '1.830.82'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").Value = '''4.993'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = '''0.6683'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '''83.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '''0.000009776'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.58%  '
$ws.Range("D17").Value = '''6.034'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '29.020.42'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").Value = '''225.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").Value = '''0.9983'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '''7.134'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").Value = '''0.9979'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").Value = '''160.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E25").Value = '  +3.78%  '
$ws.Range("D26").Value = '''8.515'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("D27").Value = '''17.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = '''1.499'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '''4.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").Value = '''0.05476'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.26%  '
$ws.Range("D32").Value = '''1.200'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").Value = '''0.7446'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").Value = '''1.137'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").Value = '''2.610'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("D37").Value = '1.241.12'
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("D38").Value = '''2.750'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("D39").Value = '''0.01784'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '''6.711'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").Value = '''0.9027'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").Value = '''0.9987'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '''101.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '1.971.32'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").Value = '''0.00000000125'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("D47").Value = '''0.5070'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("D48").Value = '''0.4053'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("D49").Value = '''0.07430'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.59%  '
$ws.Range("D50").Value = '''8.963'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '
